$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Sending cluster (A), Ligand symbol (B), Receptor symbol (C), Target cluster (D),
# then columns E..T are the numeric metrics.
# Full cross-product of sending clusters (ECs, FAPs, M2, sCs) x target clusters (FAPs, sCs).

$rows = @(
  @{ row=2;  A="ECs";  D="FAPs"; E=3; F=1;                  G=150.258513;        H=450.775539;        I=0.5395416880146598; J=0.5395416880146598; K=3; L=1;                  M=34.682839;         N=104.048517;        O=0.9919940127517238; P=0.9919940127517238; Q=5211.391814758407;  R=46902.52633282566;  S=0.5352221241405011;    T=0.5352221241405011 }
  @{ row=3;  A="ECs";  D="sCs";  E=3; F=1;                  G=150.258513;        H=450.775539;        I=0.5395416880146598; J=0.5395416880146598; K=2; L=0.6666666666666666; M=0.2799113333333333; N=0.839734;          O=0.008005987248276263; P=0.008005987248276263; Q=42.059060718514;    R=378.531546466626;   S=0.004319563874158817;  T=0.004319563874158817 }
  @{ row=4;  A="FAPs"; D="FAPs"; E=3; F=1;                  G=36.46294533333333; H=109.388836;        I=0.1309295472339256; J=0.1309295472339256; K=3; L=1;                  M=34.682839;         N=104.048517;        O=0.9919940127517238; P=0.9919940127517238; Q=1264.638462461801;  R=11381.74616215621;  S=0.1298813269483483;    T=0.1298813269483483 }
  @{ row=5;  A="FAPs"; D="sCs";  E=3; F=1;                  G=36.46294533333333; H=109.388836;        I=0.1309295472339256; J=0.1309295472339256; K=2; L=0.6666666666666666; M=0.2799113333333333; N=0.839734;          O=0.008005987248276263; P=0.008005987248276263; Q=10.20639164551378;  R=91.857524809624;    S=0.001048220285577393;  T=0.001048220285577393 }
  @{ row=6;  A="M2";   D="FAPs"; E=2; F=0.6666666666666666; G=0.2335036666666667; H=0.700511;         I=0.0008384547401380566; J=0.0008384547401380566; K=3; L=1;           M=34.682839;         N=104.048517;        O=0.9919940127517238; P=0.9919940127517238; Q=8.098570076909667;  R=72.887130692187;    S=0.0008317420821802546; T=0.0008317420821802546 }
  @{ row=7;  A="M2";   D="sCs";  E=2; F=0.6666666666666666; G=0.2335036666666667; H=0.700511;         I=0.0008384547401380566; J=0.0008384547401380566; K=2; L=0.6666666666666666; M=0.2799113333333333; N=0.839734;      O=0.008005987248276263; P=0.008005987248276263; Q=0.06536032267488889; R=0.5882429040739999; S=0.000006712657957802069; T=0.000006712657957802069 }
  @{ row=8;  A="sCs";  D="FAPs"; E=3; F=1;                  G=91.53790766666667; H=274.613723;        I=0.3286903100112765; J=0.3286903100112764; K=3; L=1;                  M=34.682839;         N=104.048517;        O=0.9919940127517238; P=0.9919940127517238; Q=3174.794513999866;  R=28573.15062599879;  S=0.3260588195806943;    T=0.3260588195806942 }
  @{ row=9;  A="sCs";  D="sCs";  E=3; F=1;                  G=91.53790766666667; H=274.613723;        I=0.3286903100112765; J=0.3286903100112764; K=2; L=0.6666666666666666; M=0.2799113333333333; N=0.839734;          O=0.008005987248276263; P=0.008005987248276263; Q=25.62249778552022;  R=230.602480069682;   S=0.002631490430582251;  T=0.002631490430582251 }
)

foreach ($r in $rows) {
  $row = $r.row
  $ws.Cells.Item($row, 1).Value = $r.A
  $ws.Cells.Item($row, 2).Value = "Timp3"
  $ws.Cells.Item($row, 3).Value = "Agtr2"
  $ws.Cells.Item($row, 4).Value = $r.D
  $ws.Cells.Item($row, 5).Value = $r.E
  $ws.Cells.Item($row, 6).Value = $r.F
  $ws.Cells.Item($row, 7).Value = $r.G
  $ws.Cells.Item($row, 8).Value = $r.H
  $ws.Cells.Item($row, 9).Value = $r.I
  $ws.Cells.Item($row, 10).Value = $r.J
  $ws.Cells.Item($row, 11).Value = $r.K
  $ws.Cells.Item($row, 12).Value = $r.L
  $ws.Cells.Item($row, 13).Value = $r.M
  $ws.Cells.Item($row, 14).Value = $r.N
  $ws.Cells.Item($row, 15).Value = $r.O
  $ws.Cells.Item($row, 16).Value = $r.P
  $ws.Cells.Item($row, 17).Value = $r.Q
  $ws.Cells.Item($row, 18).Value = $r.R
  $ws.Cells.Item($row, 19).Value = $r.S
  $ws.Cells.Item($row, 20).Value = $r.T
}
